$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.585.17"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.640.18"
$ws.Range("E3").Value = "  +1.57%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.61"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.76"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.572"
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.95"
$ws.Range("E9").Value = "  +7.86%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "3.108.64"
$ws.Range("D14").Value = "59.479.20"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.33"
$ws.Range("E15").Value = "  +3.91%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.580.15"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").Value = "  +2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "339.05"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.31"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.27"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("E24").Value = "  +2.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.29"
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("D28").Value = "0.0₃0749"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.65"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "150.96"
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.00"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.837"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "284.96"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0539"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.19"
$ws.Range("E45").Value = "  +3.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0943"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("D48").Value = "1.960.37"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.56"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.40"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.12"
$ws.Range("E51").Value = "  -0.19%  "
